# Update "想去人数" (interest count) values in column F
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 361
$ws1.Range("F3").Value = 763
$ws1.Range("F4").Value = 263
$ws1.Range("F5").Value = 784
$ws1.Range("F6").Value = 1898
$ws1.Range("F7").Value = 169

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 361
$ws4.Range("F3").Value = 763
$ws4.Range("F4").Value = 263
$ws4.Range("F7").Value = 784
$ws4.Range("F8").Value = 1898
$ws4.Range("F10").Value = 169
